# Revert "adding term 2.0 now utf-8"
# - Remove the "Include from FSIII 2" worksheet
# - Restore Version/Date/Contact/descendent-of value cells on remaining sheets

$wb = $excel.ActiveWorkbook

# Update Metadata sheet (sheet 1)
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# Update "Include from FSIII" sheet (sheet 2)
$inc = $wb.Worksheets.Item("Include from FSIII")
$inc.Range("C2").Value = "E"

# Remove the "Include from FSIII 2" worksheet entirely
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Include from FSIII 2").Delete()
$excel.DisplayAlerts = $true

# Keep the original active sheet selection (Metadata / first sheet)
$meta.Activate()
$meta.Select()
